$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 45, shifting rows 45:66 down to 46:67
$ws.Rows("45:45").Insert()

# Copy the date style (format) used by column D down into the new row's D cell
$ws.Range("D44").Copy()
$ws.Range("D45").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the new row 45 with the inserted record's data
$ws.Range("A45").Value = 5
$ws.Range("B45").Value = "Macroferia Regional de Talca"
$ws.Range("C45").Value = "Maule"
$ws.Range("D45").Value = 44529
$ws.Range("E45").Value = 7
$ws.Range("F45").Value = 100112026
$ws.Range("G45").Value = "Haba"
$ws.Range("H45").Value = "Sin especificar"
$ws.Range("I45").Value = "Primera"
$ws.Range("J45").Value = 200
$ws.Range("K45").Value = 8000
$ws.Range("L45").Value = 8000
$ws.Range("M45").Value = 8000
$ws.Range("N45").Value = "$/saco 25 kilos"
$ws.Range("O45").Value = "Región del Maule"
$ws.Range("P45").Value = 320
$ws.Range("Q45").Value = 25
$ws.Range("R45").Value = "Hortaliza"
